$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf13b"
$ws.Range("C2").Value = "Tnfrsf13c"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 2.606130666666667
$ws.Range("H2").Value = 7.818392
$ws.Range("I2").Value = 0.2943337015143609
$ws.Range("J2").Value = 0.2943337015143609
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.091090666666667
$ws.Range("N2").Value = 3.273272
$ws.Range("O2").Value = 0.9548729479813977
$ws.Range("P2").Value = 0.9548729479813977
$ws.Range("Q2").Value = 2.843524846513778
$ws.Range("R2").Value = 25.591723618624
$ws.Range("S2").Value = 0.2810512892552945
$ws.Range("T2").Value = 0.2810512892552946

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfsf13b"
$ws.Range("C3").Value = "Tnfrsf13c"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 2.606130666666667
$ws.Range("H3").Value = 7.818392
$ws.Range("I3").Value = 0.2943337015143609
$ws.Range("J3").Value = 0.2943337015143609
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05156466666666667
$ws.Range("N3").Value = 0.154694
$ws.Range("O3").Value = 0.04512705201860229
$ws.Range("P3").Value = 0.04512705201860228
$ws.Range("Q3").Value = 0.1343842591164444
$ws.Range("R3").Value = 1.209458332048
$ws.Range("S3").Value = 0.01328241225906632
$ws.Range("T3").Value = 0.01328241225906632

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfsf13b"
$ws.Range("C4").Value = "Tnfrsf13c"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.019967333333334
$ws.Range("H4").Value = 18.059902
$ws.Range("I4").Value = 0.6798888831164529
$ws.Range("J4").Value = 0.6798888831164529
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.091090666666667
$ws.Range("N4").Value = 3.273272
$ws.Range("O4").Value = 0.9548729479813977
$ws.Range("P4").Value = 0.9548729479813977
$ws.Range("Q4").Value = 6.568330171038222
$ws.Range("R4").Value = 59.114971539344
$ws.Range("S4").Value = 0.6492075021211873
$ws.Range("T4").Value = 0.6492075021211873

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfsf13b"
$ws.Range("C5").Value = "Tnfrsf13c"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.019967333333334
$ws.Range("H5").Value = 18.059902
$ws.Range("I5").Value = 0.6798888831164529
$ws.Range("J5").Value = 0.6798888831164529
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05156466666666667
$ws.Range("N5").Value = 0.154694
$ws.Range("O5").Value = 0.04512705201860229
$ws.Range("P5").Value = 0.04512705201860228
$ws.Range("Q5").Value = 0.3104176088875556
$ws.Range("R5").Value = 2.793758479988
$ws.Range("S5").Value = 0.03068138099526558
$ws.Range("T5").Value = 0.03068138099526558

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Tnfsf13b"
$ws.Range("C6").Value = "Tnfrsf13c"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.228242
$ws.Range("H6").Value = 0.6847260000000001
$ws.Range("I6").Value = 0.02577741536918619
$ws.Range("J6").Value = 0.02577741536918619
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.091090666666667
$ws.Range("N6").Value = 3.273272
$ws.Range("O6").Value = 0.9548729479813977
$ws.Range("P6").Value = 0.9548729479813977
$ws.Range("Q6").Value = 0.2490327159413334
$ws.Range("R6").Value = 2.241294443472
$ws.Range("S6").Value = 0.0246141566049158
$ws.Range("T6").Value = 0.0246141566049158

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Tnfsf13b"
$ws.Range("C7").Value = "Tnfrsf13c"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.228242
$ws.Range("H7").Value = 0.6847260000000001
$ws.Range("I7").Value = 0.02577741536918619
$ws.Range("J7").Value = 0.02577741536918619
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.05156466666666667
$ws.Range("N7").Value = 0.154694
$ws.Range("O7").Value = 0.04512705201860229
$ws.Range("P7").Value = 0.04512705201860228
$ws.Range("Q7").Value = 0.01176922264933333
$ws.Range("R7").Value = 0.105923003844
$ws.Range("S7").Value = 0.001163258764270383
$ws.Range("T7").Value = 0.001163258764270383

